$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, copying the header style/format from G1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill in the Save column values for rows 2-14
$values = @(1, 0, 0, 0, 0, 1, 1, 0, 1, 0, 0, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
